$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: restyle from the "4/5" style family to the "10/11" style family ---
# Copy cell formats from row 6 (which already uses styles 10/11) onto row 13,
# leaving the existing values untouched.
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial(-4122) | Out-Null

# --- New shared strings must be created in this exact order so the shared
#     string table indices line up with the target workbook ---
$ws.Range("A14").Value = "SCRIPT/P01P04A/um2105.ssb"
$ws.Range("C14").Value = " You were out on some kind of\njourney, weren\'t you?"
$ws.Range("C15").Value = " I just got back from a long\njourney myself."
$ws.Range("D15").Value = " Я тоже не так давно вернулся из\nзатяжного путешествия."
$ws.Range("D14").Value = " Вы куда-то путешествовали, да?"
$ws.Range("E14").Value = " Âú ëôäà-óï ðôóåšåòóâïâàìé, äà?"
$ws.Range("E15").Value = " Ÿ óïçå îå óàë äàâîï âåñîôìòÿ éè\nèàóÿçîïãï ðôóåšåòóâéÿ."
$ws.Range("C16").Value = " I wonder where I\'ll go\nexploring next…"
$ws.Range("A16").Value = "SCRIPT/P01P04A/um2205.ssb"
$ws.Range("D16").Value = " Интересно, куда я отправлюсь\nна вылазку в следующий раз..."
$ws.Range("E16").Value = " Éîóåñåòîï, ëôäà ÿ ïóðñàâìýòû\nîà âúìàèëô â òìåäôýþéê ñàè..."

# --- numeric "line number" column ---
$ws.Range("B14").Value = 217
$ws.Range("B15").Value = 220
$ws.Range("B16").Value = 198

# --- formatting: rows 14 & 16 mirror row 8's style family (4/5), row 15
#     mirrors row 9's style family (8/9) ---
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A14:E14").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null

$ws.Range("A9:E9").Copy() | Out-Null
$ws.Range("A15:E15").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(14).RowHeight = 43.2
$ws.Rows.Item(15).RowHeight = 21.6
$ws.Rows.Item(16).RowHeight = 26.4

# --- sheet view adjustments ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D16").Select() | Out-Null
